$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DM_Stat (C) and P_Value (D) columns for rows 2-11
$ws.Range("C2").Value = -0.9310486228775472
$ws.Range("D2").Value = 0.3619335799283108

$ws.Range("C3").Value = 0.138406833522296
$ws.Range("D3").Value = 0.8911779709330114

$ws.Range("C4").Value = 0.89393404185339
$ws.Range("D4").Value = 0.3810327363097101

$ws.Range("C5").Value = 0.2704347423026691
$ws.Range("D5").Value = 0.7893466636286499

$ws.Range("C6").Value = 0.7209216754196497
$ws.Range("D6").Value = 0.478549207354106

$ws.Range("C7").Value = 2.197211383216863
$ws.Range("D7").Value = 0.03883072318088576
$ws.Range("G7").Value = "Sí"

$ws.Range("C8").Value = 1.021830866778255
$ws.Range("D8").Value = 0.3179666917529407

$ws.Range("C9").Value = 0.5570140928567286
$ws.Range("D9").Value = 0.5831407916316635

$ws.Range("C10").Value = 0.1678486005809995
$ws.Range("D10").Value = 0.868235815161116

$ws.Range("C11").Value = -0.5698919051348892
$ws.Range("D11").Value = 0.5745251345366456
